$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.858.36"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.71%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.906.47"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.26%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.9996"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.50%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'313.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.76%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.9994"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.45%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4995"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +3.79%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3811"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.03%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.07273"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.22%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.9113"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -2.52%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'21.01"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.83%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = "'1.941.61"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +1.11%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").Value = "'0.07636"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.82%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'5.494"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.04%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'91.97"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.25%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.9996"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.53%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.000008732"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.50%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.9992"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.44%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'27.892.29"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.69%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -1.15%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'5.180"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.31%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'10.85"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.71%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'6.565"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.40%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'152.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -1.78%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'1.860"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -2.98%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'2.224"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +3.58%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'18.40"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.59%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'115.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -1.39%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'4.913"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.18%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.09013"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.67%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'3.200"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -3.25%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'4.829"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +3.22%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -2.76%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.7741"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.58%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.02086"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +1.41%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("B36").Value = 'Frax'
$ws.Range("C36").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D36").Value = "'0.9986"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.64%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").Value = "'2.556"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -2.84%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'3.056"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +1.73%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D39").Value = "'1.094"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -1.66%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").Value = "'0.5553"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +1.00%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = 'Hedera'
$ws.Range("C41").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D41").Value = "'0.05273"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.91%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = "'6.894"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -2.08%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = 'Aptos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D43").Value = "'8.492"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.02%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = 'Algorand'
$ws.Range("C44").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D44").Value = "'0.1521"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.64%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").Value = "'112.66"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +3.87%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.4835"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.06%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = "'10.59"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.43%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = 'PaxDollar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D48").Value = "'0.9992"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.47%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").Value = "'1.634"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.99%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").Value = "'67.44"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.28%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = "'0.06053"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.42%  "
$ws.Range("E51").Style = "Normal"
